$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.899.18'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.628.16'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.520'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.46'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.07%  '
$ws.Range("E9").Value = '  +3.46%  '
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").Value = '1.862.72'
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("D13").Value = '1.630.56'
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("E14").Value = '  +6.23%  '
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("D16").Value = '29.942.37'
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '9.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +19.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  +3.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.82%  '
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.38%  '
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("E28").Value = '  +3.21%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("E31").Value = '  +6.48%  '
$ws.Range("E32").Value = '  +3.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").Value = '1.428.99'
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("E35").Value = '  +6.46%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("E39").Value = '  +3.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.556'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.06%  '
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("E42").Value = '  +3.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.92%  '
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '69.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.28%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("E48").Value = '  +2.55%  '
$ws.Range("D49").Value = '1.770.11'
$ws.Range("E49").Value = '  +1.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '89.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.98%  '
$ws.Range("D51").Value = '0.0₆0107'
$ws.Range("E51").Value = '  +1.24%  '
